# Datebook_jph.xlsx edit: add 0706 comment, 0707 goal
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13 (2021-07-06): fill in Time/Comment, update Goal text ---
$ws.Range("C13").Value = "08:32`n18:03"
$ws.Range("C13").WrapText = $true

$ws.Range("D13").Value = "1. 추가 자료조사(관련 코드)`n2. FlowChart 수정 및 보완"
$ws.Range("D13").WrapText = $true

$ws.Range("E13").Value = "1. 필요 자료 추가 조사`n2. 한계점 및 필요기술 추가`n - 차량 블루투스와 어플 연결 어떻게?`n - 연결이 된다면 원하는 정보만 어떻게 가져 올건지?"
$ws.Range("E13").WrapText = $true

# --- Row 14 (2021-07-07): fill in Goal text ---
$ws.Range("D14").Value = "1. 앱 화면별 Sequence 작성`n2. 앱 화면별 어떻게 구성할 것인지 구상"
$ws.Range("D14").WrapText = $true

# --- Row heights to match the wrapped content ---
$ws.Rows.Item(13).RowHeight = 52.700000000000003
$ws.Rows.Item(14).RowHeight = 26.350000000000001

# --- View state: scroll position + active selection ---
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("E14").Select()
